# Applies the edit described by the diff:
#  - 4 new shared strings were inserted ("Gaussian-Quadrature" relocated + 3 new
#    "Spiral-..." labels) right after "Ring Perpendicular to TD" in the shared
#    string table, and the old trailing "Gaussian-Quadrature" entry was removed.
#    Net effect on the worksheet: rows that used to read
#       NoRotation-tilt60deg / Rotation-NoTilt / Rotation-60detTilt /
#       HexGrid-90degTilt5degRes / HexGrid-90degTilt22p5degRes /
#       HexGrid-60degTilt5degRes / Gaussian-Quadrature (rows 10-16)
#    now read
#       Gaussian-Quadrature / Spiral-90deg-10rot-5space / Spiral-90deg-15rot-5space /
#       Spiral-90deg-10rot-3space / NoRotation-tilt60deg / Rotation-NoTilt /
#       Rotation-60detTilt (rows 10-16),
#    and three brand new rows (17-19) are appended for the HexGrid entries,
#    each with all-1 data across columns C:M (same pattern as every other row).
#  - The sheet dimension grows from A1:M16 to A1:M19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label ordering for rows 10-16 (column B), reflecting the updated
# shared-string table ordering.
$labels = @(
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt"
)

$row = 10
foreach ($label in $labels) {
    $ws.Cells.Item($row, 2).Value = $label
    $row++
}

# Append the three new rows (17, 18, 19) for the HexGrid entries that used
# to occupy rows 13-15, each populated with the same all-ones pattern used
# throughout the table.
$newRows = @(
    @{ Index = 15; Label = "HexGrid-90degTilt5degRes" },
    @{ Index = 16; Label = "HexGrid-90degTilt22p5degRes" },
    @{ Index = 17; Label = "HexGrid-60degTilt5degRes" }
)

$row = 17
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry.Index
    # Copy the column-A cell formatting (bold, centered, bordered) from the
    # row above so the new rows match the rest of the table's A-column style.
    $ws.Range("A" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($row, 2).Value = $entry.Label

    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }

    $row++
}

$excel.CutCopyMode = $false
